$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Historical Invoice Template")
$ws1.Range("B2").Value = "Invoice1338448"
$ws1.Range("C2").Value = "Invoice1302475"
$ws1.Range("E2").Value = "Lassie9p48"
$ws1.Range("F2").Value = "Chloeb3m4"

$ws2 = $wb.Worksheets.Item("Historical PO Template")
$ws2.Range("B2").Value = "Invoice1302475"
$ws2.Range("C2").Value = "Invoice1338448"
$ws2.Range("E2").Value = "Lassie9p48"
$ws2.Range("F2").Value = "Chloeb3m4"
